$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row.
# The update bumps this date by one day (serial 46060 -> 46061, i.e. 2026-02-07 -> 2026-02-08)
# for every row in the table, from row 2 through row 550.
$lastRow = 550
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46061
